# feat: add 2022-Q3 data
#
# 1) Insert a brand-new worksheet "2022-Q3" right before the existing
#    "2022-Q2" sheet (which pushes it and every later quarter sheet one
#    slot to the right / up the sheetId order).
# 2) Fill the new sheet with the quarterly fund-holding table.
# 3) Insert a new row 2 into the "总计" (totals) summary sheet and fill
#    it with the 2022-Q3 totals, shifting the older rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: insert the new "2022-Q3" worksheet before "2022-Q2"
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($q2Sheet)
$newSheet.Name = "2022-Q3"

# Match the sheetPr/pageMargins boilerplate every other quarter sheet has.
$newSheet.Outline.SummaryBelow = $true
$newSheet.Outline.SummaryRight = $true
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# Step 2: populate the new sheet
# ---------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $newSheet.Cells.Item(1, $i + 2)  # headers start at column B
    $cell.Value = $headers[$i]
}
$newSheet.Range("B1:H1").Font.Bold = $true
$newSheet.Range("B1:H1").HorizontalAlignment = -4108

$dataRows = @(
    @(0, "006218", "富国生物医药科技混合A", "7.10", "89.32", "5.88", "0.4175", 4),
    @(1, "100016", "富国天源沪港深平衡混合A", "4.99", "70.11", "3.18", "0.1587", 4),
    @(2, "011308", "富国生物医药科技混合C", "1.52", "89.32", "5.88", "0.0894", 4),
    @(3, "005108", "圆信永丰双利优选定期开放灵活配置混合", "0.61", "91.10", "4.09", "0.0249", 5),
    @(4, "001965", "圆信永丰兴源灵活配置混合A", "0.48", "93.86", "4.28", "0.0205", 5),
    @(5, "001966", "圆信永丰兴源灵活配置混合C", "0.26", "93.86", "4.28", "0.0111", 5),
    @(6, "006274", "圆信永丰医药健康混合", "0.16", "93.60", "4.15", "0.0066", 5),
    @(7, "001563", "华富健康文娱灵活配置混合", "0.13", "93.80", "2.83", "0.0037", 10),
    @(8, "015655", "富荣医药健康混合A", "0.13", "82.11", "1.62", "0.0021", 8),
    @(9, "015656", "富荣医药健康混合C", "0.01", "82.11", "1.62", "0.0002", 8),
    @(10, "014931", "富国天源沪港深平衡混合C", "0.00", "70.11", "3.18", 0, 4)
)

$rowIndex = 2
foreach ($row in $dataRows) {
    $colIndex = 1
    foreach ($value in $row) {
        $cell = $newSheet.Cells.Item($rowIndex, $colIndex)
        # Numeric-looking strings (fund codes with leading zeros, ratios,
        # sizes, ...) must stay text - otherwise Excel's normal typed-value
        # coercion would turn "006218" into 6218 or "7.10" into 7.1.
        if ($value -is [string]) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $value
        $colIndex++
    }
    $rowIndex++
}

# Column A (row index) uses the same bold/centered "index" look as the
# header row and as every other quarter sheet's index column.
$indexRange = $newSheet.Range("A2:A12")
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160

$newSheet.Range("A1").Select()

# ---------------------------------------------------------------------
# Step 3: update the "总计" summary sheet
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 11
$totalSheet.Cells.Item(2, 4).Value = 0.73

# A2 (new row's index cell) gets the same bold/centered look the index
# column uses on every other row of this sheet.
$totalA2 = $totalSheet.Cells.Item(2, 1)
$totalA2.Font.Bold = $true
$totalA2.HorizontalAlignment = -4108
$totalA2.VerticalAlignment = -4160

# Column A is a plain positional index (0, 1, 2, ...), not a value that
# travels with its row - renumber every row below the new one so it
# stays a contiguous 0-based sequence.
for ($r = 3; $r -le 8; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
